$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data range (A2:B53) by column A, ascending, stable - matches
# the manual "Data > Sort" operation the author performed.
$ws.Range("A2:B53").Sort($ws.Range("A2:A53"))

# Restore the cell selection to where the author ended up after sorting.
$ws.Range("D47").Select() | Out-Null
